# edit.ps1 - applies the "rettelser i resultater og diskussion" changes
# to the active Word document via COM-interop (Range.InsertXML).
#
# Strategy: for every paragraph whose internal run/proofErr structure
# changes in the target diff, replace that paragraph's Range content
# wholesale with freshly authored OOXML that reproduces the exact
# <w:r>/<w:t>/<w:proofErr> structure requested by the diff. This keeps
# every other paragraph (and all drawings/fields elsewhere in the
# document) completely untouched.

$d = $word.ActiveDocument

function Set-RangeXml($range, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        $bodyXml +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- Paragraph 2: "Igennem intergrationstesten ..." -------------------
# Fix the "intergrationstesten" -> "integrationstesten" typo. The
# author's correction splits the word across a run boundary
# ("Igennem inte" | "grationstesten ...") and removes the spell-check
# proofErr wrapper that used to flag the misspelling.
$para2 = $d.Paragraphs(2).Range
$body2 = @'
<w:body><w:p>
  <w:r><w:t>Igennem inte</w:t></w:r>
  <w:r><w:t xml:space="preserve">grationstesten og accepttesten er der opn&#229;et resultater for projektet. Igennem accepttesten er der opstillet krav til hvordan blodtryksm&#229;ler systemet skal opf&#248;re sig i forhold til </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>use</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> cases, og hvordan dette kommer til udtryk visuelt. Resultaterne for projektet er visuelle resultater af accepttesten, og ses i dette afsnit som </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>screendumps</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">. </w:t></w:r>
</w:p></w:body>
'@
Set-RangeXml $para2 $body2

# --- Paragraph 3: "I acceptesten er det f\u00f8rste step ..." ------------
# Fix the "acceptesten" -> "accepttesten" typo (missing double-t),
# splitting it as "I accept" | "t" | "esten ..." and removing the
# spell-check proofErr wrapper.
$para3 = $d.Paragraphs(3).Range
$body3 = @'
<w:body><w:p>
  <w:r><w:t>I accept</w:t></w:r>
  <w:r><w:t>t</w:t></w:r>
  <w:r><w:t xml:space="preserve">esten er det f&#248;rste step for accepttest af </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>use</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> case 1, er at v&#230;lge en v&#230;rdi p&#229; vands&#248;jlen og kalibrer efter den. Herved skal den afl&#230;ste sp&#230;nding og trykket i vands&#248;jlen indtastes. N&#229;r programmet starter, vises start sk&#230;rmen neden for.</w:t></w:r>
</w:p></w:body>
'@
Set-RangeXml $para3 $body3

# --- Paragraph 25: "Fra de forg\u00e5ende figur ..." ---------------------
# "kravene" -> "accepttestens krav", split across 3 runs.
$para25 = $d.Paragraphs(25).Range
$body25 = @'
<w:body><w:p>
  <w:r><w:t>Fra de forg&#229;ende figur kan det konkluderes at blodtryksm&#229;l</w:t></w:r>
  <w:r><w:t xml:space="preserve">ersystemet lever op til accepttestens krav </w:t></w:r>
  <w:r><w:t xml:space="preserve">omkring, at kunne vise et blodtrykssignalet kontinuert, samt s&#230;tte et digitalt filter til og fra. </w:t></w:r>
</w:p></w:body>
'@
Set-RangeXml $para25 $body25

# --- Paragraphs 26-30: reworked / reordered discussion --------------
# Paragraph 26 ("Blodtryksm\u00e5lersystemet lever ogs\u00e5 ...") keeps its
# opening but the closing sentence about nulpunktsv\u00e6rdi is reworded
# and split into several runs.
# Paragraph 27 ("Der kan ogs\u00e5 \u00e6ndres gr\u00e6nsev\u00e6rdierne ...") gets its
# "blodtryksm\u00e5lersystetemet" typo fixed (spell-check proofErr removed)
# and is split across 3 runs.
# Paragraph 28 ("Der kan ogs\u00e5 uds\u00e6ttes alarmen ...") keeps its text but
# loses the trailing _GoBack bookmark.
# Paragraph 29 ("En anden ting ...") is unchanged, but now sits before
# the brand new paragraph 30 (instead of before a trailing empty
# paragraph).
# Paragraph 30 is an entirely new paragraph ("Det kan diskuteres om
# brugeren skal have lov til at kalibrere p\u00e5 startsk\u00e6rmen ...") which
# now also hosts the relocated _GoBack bookmark; the old trailing
# empty paragraph is removed.
$startRange = $d.Paragraphs(26).Range
$endRange = $d.Paragraphs(30).Range
$fullRange = $d.Range($startRange.Start, $endRange.End)
$bodyTail = @'
<w:body>
<w:p>
  <w:r><w:t xml:space="preserve">Blodtryksm&#229;lersystemet lever ogs&#229; op til kravene omkring at kunne kalibrer systemet, og dette bliver gjort p&#229; startsk&#230;rmen, samt nulpunktsjuster, som sker p&#229; </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>HovedGUI</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>. Det kan diskuteres om det var smartere at systemet selv skal kunne indl&#230;se v&#230;rdien for nulpunktjustering, og derefter bare gange denne v&#230;rdi p&#229; blodtrykssignalet, n&#229;r der trykkes p&#229; nulpunkts justerings</w:t></w:r>
  <w:r><w:t xml:space="preserve"> knappen, end at man selv skal afl&#230;se v&#230;rdien fra et andet program. Dog lever nulpunks justerings knappen, op til kravet omkring at nulpunkts juster blodtryksm&#229;ler systemet, da alle blodtrykssignalv&#230;rdierne</w:t></w:r>
  <w:r><w:t xml:space="preserve"> bliver</w:t></w:r>
  <w:r><w:t xml:space="preserve"> nulpunksjusteret efter</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>den</w:t></w:r>
  <w:r><w:t xml:space="preserve"> indtastet nulpunktsv&#230;rdi</w:t></w:r>
  <w:r><w:t xml:space="preserve">. </w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>Der kan ogs&#229; &#230;ndres gr&#230;nsev&#230;r</w:t></w:r>
  <w:r><w:t>dierne for blodtryksm&#229;lersyste</w:t></w:r>
  <w:r><w:t xml:space="preserve">met, herved er dette krav opfyldt, og der kommer en alarm n&#229;r v&#230;rdierne for enten systolen eller diastolen overskrider de satte gr&#230;nsev&#230;rdier. </w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">Der kan ogs&#229; uds&#230;ttes alarmen for blodtryksm&#229;leren i et minut, ved at trykke p&#229; alarm knappen oppe i h&#248;jre hj&#248;rne af </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>HovedGUI&#8217;en</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">. </w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">En anden ting som blodtryksm&#229;lersystemet ogs&#229; kan er at vise en timer, som starter n&#229;r man har trykket p&#229; start knappen, og stopper n&#229;r der trykkes p&#229; sluk knappen. Det kan diskuteres hvor smart det er, om timeren skal nulstille, n&#229;r der trykkes p&#229; start knappen igen. Der er valgt i dette projekt at timeren starter fra det stoppet tidspunkt, n&#229;r t&#230;nd bliver trykket igen, dette er valgt fordi, at det giver bedst mening at man skal kunne starte fra hvor man slap. Dog giver den funktion med at starte timeren igen, ikke s&#229; meget mening ude i den virkelige verden, da an&#230;stesi sygeplejerskerne ikke starter m&#229;lingen igen, efter at have stoppet m&#229;lingen. </w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">Det kan diskuteres om brugeren skal have lov til at kalibrer p&#229; startsk&#230;rmen, da kalibrering skal fortages af </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>en</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> servicepersonale. Herved skulle der kun have v&#230;ret log ind p&#229; startsk&#230;rmen, og have lavet et servicevindue som kun kan betjenes af servicepersonalet. </w:t></w:r>
  <w:r><w:t xml:space="preserve">Der er taget h&#248;jde for denne problemstilling i koden, ved at have lavet en </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>config</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> fil, hvor der kan &#230;ndres kalibreringstallet i. Det betyder, at hvis det </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>sundhedsfagligt</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> personale springer over at kalibrer p&#229; startsk&#230;rmen, kalibrer systemet automatisk efter kalibreringstallet i </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>config</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> filen. Meningen med </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>config</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> filen, er at det kun er servicepersonalet der skal kunne tilg&#229; denne fil, </w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:t xml:space="preserve">og kunne &#230;ndre kalibreringstallet. </w:t></w:r>
</w:p>
</w:body>
'@
Set-RangeXml $fullRange $bodyTail
